$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; this shifts existing rows 14-63 down to 15-64
$ws.Rows.Item(14).Insert()

# Fill in the new row 14 with the new data record
$newDate = Get-Date -Year 2022 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = $newDate
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = 100112021
$ws.Cells.Item(14, 7).Value = "Ají"
$ws.Cells.Item(14, 8).Value = "Americana (o)"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 9500
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = 9750
$ws.Cells.Item(14, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 650
$ws.Cells.Item(14, 17).Value = 15
$ws.Cells.Item(14, 18).Value = "Hortaliza"
